# Updates crypto price/volume data to reflect the latest scrape.
# Rows 47 and 48 also swap coin identity (ordi <-> Aave reorder by rank).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.154.80"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "2.305.27"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  -0.19%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.05"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -0.66%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.40"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +2.91%  "
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("E8").Value = "  -0.16%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.614"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +1.12%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.10"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +1.50%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0916"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +0.23%  "
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.43"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  +2.08%  "
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +1.08%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.984"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +1.38%  "
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.50"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "2.655.26"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").Value = "2.307.35"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "42.146.31"
$ws.Range("E18").Value = "  -0.75%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.74"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("E20").Value = "  +0.78%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.20"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -4.13%  "
$ws.Range("E22").Value = "  -0.16%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "262.97"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +1.55%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.34"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +2.67%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.97"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +4.93%  "
$ws.Range("E26").Value = "  +0.42%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.03"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -1.69%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.32"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +4.96%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.95"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +0.08%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.74"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +4.55%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.69"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("E32").Value = "  +2.09%  "
$ws.Range("E33").Value = "  -1.78%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.92"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +0.06%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.121"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +10.72%  "
$ws.Range("E36").Value = "  +0.45%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.65"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +3.25%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.98"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +13.99%  "
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("E40").Value = "  -1.22%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.38"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +19.85%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.50"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +3.59%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.02"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +4.96%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.229"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("E45").Value = "  -0.03%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.51"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +8.70%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "114.59"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "80.64"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +13.93%  "
$ws.Range("E49").Value = "  +2.32%  "
$ws.Range("E50").Value = "  -1.65%  "
